$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.477.19"
$ws.Range("E2").Value = "  +3.05%  "
$ws.Range("D3").Value = "2.543.78"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.06"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.98"
$ws.Range("E6").Value = "  +3.44%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("D9").Value = "2.543.77"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  +2.62%  "
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.349"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.09"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "3.010.81"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "68.438.58"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "2.549.47"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.03"
$ws.Range("E19").Value = "  +3.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.57"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.76"
$ws.Range("E21").Value = "  +6.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.74"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.34"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.04"
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").Value = "2.682.46"
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  +2.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "544.21"
$ws.Range("E31").Value = "  +2.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.30"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("E33").Value = "  +2.54%  "
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.59"
$ws.Range("E38").Value = "  +0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.95"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.70"
$ws.Range("E40").Value = "  +1.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.357"
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.566"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.29"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0283"
$ws.Range("E48").Value = "  +3.80%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  +1.00%  "
